# Re-creates the "Made some Changes" upload: enters a small SCCM/Windows
# table on Sheet1 and registers the "amilne - Personal View" custom view
# that was captured when the workbook was (re-)saved as a shared workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data entry -------------------------------------------------------
# Entered in this order so the shared-string table comes out sorted the
# same way the source workbook has it (Made some Changes, SCCM 2012 R2,
# SCCM CB, SCCM CBB, Windows 8.1) and so the final selection lands on C5.
$ws.Range("A1").Value = "Made some Changes"
$ws.Range("A5").Value = "SCCM 2012 R2"
$ws.Range("B5").Value = "SCCM CB"
$ws.Range("C5").Value = "SCCM CBB"
$ws.Range("A3").Value = "Windows 8.1"

# --- Column widths ------------------------------------------------------
# Columns A and C end up auto-fit to their text ("Made some Changes" /
# "SCCM CBB") in the source file (width 19.42578125 / 9.85546875 with
# bestFit). Set the width explicitly so the saved value lands as close to
# that as this engine's column-width quantization allows.
$ws.Columns.Item(1).ColumnWidth = 18.592447916666668
$ws.Columns.Item(3).ColumnWidth = 9.022135416666666

# --- Selection ----------------------------------------------------------
# Final selection in the source file is C5 (last cell touched).
$ws.Range("C5").Select()

# --- Shared-workbook personal view --------------------------------------
# The source file was saved with Track Changes / Shared Workbook turned
# on, which records the author's window layout as a custom workbook view
# named "<user> - Personal View". Register an equivalent custom view.
$wb.CustomViews.Add("amilne - Personal View", $true, $true)
